# "cambio en hoja de resultados"
# Update the numeric results in column B (counts per Zona/Region/Plaza/Permiso row)
# to reflect the new figures, and move the sheet's viewport/selection down to the
# last block of rows (B23:B25), matching the author's final view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B4").Value = 3

$ws.Range("B11").Value = 4
$ws.Range("B12").Value = 6

$ws.Range("B15").Value = 8
$ws.Range("B16").Value = 8
$ws.Range("B17").Value = 8

$ws.Range("B18").Value = 5
$ws.Range("B19").Value = 5
$ws.Range("B20").Value = 5

$ws.Range("B23").Value = 6
$ws.Range("B24").Value = 6
$ws.Range("B25").Value = 6

# Scroll the window so row 10 is at the top and select B23:B25 (activeCell = B23),
# matching the updated sheetView/selection in the saved workbook.
$ws.Activate()
$win = $wb.Windows.Item(1)
$win.ScrollRow = 10
$win.ScrollColumn = 1
$ws.Range("B23:B25").Select()
